# EIA Table 1.12.B refresh: Year-to-Date period rolls from October to November
# (2016 & 2015 YTD), 2017-01-31 EPM run -- updates headline dates, column headers,
# and the revised YTD generation / percentage-change figures for each affected row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table_1_12_B")

# --- Title / column header text -------------------------------------------------
$ws.Range("A2").Value = "by State, by Sector, Year-to-Date through November 2016 and 2015 (Thousand Megawatthours)"

$ws.Range("B6").Value = "November 2016 YTD"
$ws.Range("C6").Value = "November 2015 YTD"
$ws.Range("E6").Value = "November 2016 YTD"
$ws.Range("F6").Value = "November 2015 YTD"
$ws.Range("G6").Value = "November 2016 YTD"
$ws.Range("H6").Value = "November 2015 YTD"
$ws.Range("I6").Value = "November 2016 YTD"
$ws.Range("J6").Value = "November 2015 YTD"
$ws.Range("K6").Value = "November 2016 YTD"
$ws.Range("L6").Value = "November 2015 YTD"

# --- Revised YTD data (Thousand Megawatthours) and percentage-change figures -----
# Row 7: New England
$ws.Range("B7").Value = -453
$ws.Range("C7").Value = -441
$ws.Range("D7").Value = 0.027
$ws.Range("G7").Value = -453
$ws.Range("H7").Value = -441

# Row 8: Connecticut
$ws.Range("B8").Value = 3
$ws.Range("D8").Value = -1.748
$ws.Range("G8").Value = 3

# Row 10: Massachusetts
$ws.Range("B10").Value = -456
$ws.Range("C10").Value = -436
$ws.Range("D10").Value = 0.044
$ws.Range("G10").Value = -456
$ws.Range("H10").Value = -436

# Row 14: Middle Atlantic
$ws.Range("B14").Value = -1121
$ws.Range("C14").Value = -1025
$ws.Range("D14").Value = 0.094
$ws.Range("E14").Value = -612
$ws.Range("F14").Value = -550
$ws.Range("G14").Value = -508
$ws.Range("H14").Value = -475

# Row 15: New Jersey
$ws.Range("B15").Value = -189
$ws.Range("C15").Value = -160
$ws.Range("D15").Value = 0.178
$ws.Range("E15").Value = -189
$ws.Range("F15").Value = -160

# Row 16: New York
$ws.Range("B16").Value = -424
$ws.Range("C16").Value = -389
$ws.Range("E16").Value = -424
$ws.Range("F16").Value = -389

# Row 17: Pennsylvania
$ws.Range("B17").Value = -508
$ws.Range("C17").Value = -475
$ws.Range("D17").Value = 0.071
$ws.Range("G17").Value = -508
$ws.Range("H17").Value = -475

# Row 18: East North Central
$ws.Range("B18").Value = -694
$ws.Range("C18").Value = -449
$ws.Range("D18").Value = 0.545
$ws.Range("E18").Value = -694
$ws.Range("F18").Value = -449

# Row 21: Michigan
$ws.Range("B21").Value = -694
$ws.Range("C21").Value = -449
$ws.Range("D21").Value = 0.545
$ws.Range("E21").Value = -694
$ws.Range("F21").Value = -449

# Row 24: West North Central
$ws.Range("B24").Value = 195
$ws.Range("C24").Value = 271
$ws.Range("D24").Value = -0.278
$ws.Range("E24").Value = 195
$ws.Range("F24").Value = 271

# Row 28: Missouri
$ws.Range("B28").Value = 195
$ws.Range("C28").Value = 271
$ws.Range("D28").Value = -0.278
$ws.Range("E28").Value = 195
$ws.Range("F28").Value = 271

# Row 32: South Atlantic
$ws.Range("B32").Value = -2859
$ws.Range("C32").Value = -2626
$ws.Range("D32").Value = 0.088
$ws.Range("E32").Value = -2859
$ws.Range("F32").Value = -2626

# Row 36: Georgia
$ws.Range("B36").Value = -851
$ws.Range("C36").Value = -826
$ws.Range("D36").Value = 0.031
$ws.Range("E36").Value = -851
$ws.Range("F36").Value = -826

# Row 39: South Carolina
$ws.Range("B39").Value = -913
$ws.Range("C39").Value = -848
$ws.Range("D39").Value = 0.077
$ws.Range("E39").Value = -913
$ws.Range("F39").Value = -848

# Row 40: Virginia
$ws.Range("B40").Value = -1094
$ws.Range("C40").Value = -953
$ws.Range("D40").Value = 0.148
$ws.Range("E40").Value = -1094
$ws.Range("F40").Value = -953

# Row 42: East South Central
$ws.Range("B42").Value = -662
$ws.Range("C42").Value = -483
$ws.Range("D42").Value = 0.371
$ws.Range("E42").Value = -662
$ws.Range("F42").Value = -483

# Row 46: Tennessee
$ws.Range("B46").Value = -662
$ws.Range("C46").Value = -483
$ws.Range("D46").Value = 0.371
$ws.Range("E46").Value = -662
$ws.Range("F46").Value = -483

# Row 47: West South Central
$ws.Range("B47").Value = -40
$ws.Range("C47").Value = -41
$ws.Range("D47").Value = -0.015
$ws.Range("E47").Value = -40
$ws.Range("F47").Value = -41

# Row 48: Arkansas
$ws.Range("B48").Value = 39
$ws.Range("D48").Value = 0.33
$ws.Range("E48").Value = 39

# Row 50: Oklahoma
$ws.Range("B50").Value = -79
$ws.Range("C50").Value = -70
$ws.Range("D50").Value = 0.128
$ws.Range("E50").Value = -79
$ws.Range("F50").Value = -70

# Row 52: Mountain
$ws.Range("B52").Value = -202
$ws.Range("C52").Value = -181
$ws.Range("D52").Value = 0.116
$ws.Range("E52").Value = -202
$ws.Range("F52").Value = -181

# Row 53: Arizona
$ws.Range("B53").Value = 68
$ws.Range("C53").Value = 75
$ws.Range("D53").Value = -0.09
$ws.Range("E53").Value = 68
$ws.Range("F53").Value = 75

# Row 54: Colorado
$ws.Range("B54").Value = -270
$ws.Range("C54").Value = -256
$ws.Range("E54").Value = -270
$ws.Range("F54").Value = -256

# Row 61: Pacific Contiguous
$ws.Range("B61").Value = -99
$ws.Range("C61").Value = 163
$ws.Range("D61").Value = -1.606
$ws.Range("E61").Value = -99
$ws.Range("F61").Value = 163

# Row 62: California
$ws.Range("B62").Value = -97
$ws.Range("C62").Value = 126
$ws.Range("D62").Value = -1.775
$ws.Range("E62").Value = -97
$ws.Range("F62").Value = 126

# Row 64: Washington
$ws.Range("B64").Value = -2
$ws.Range("C64").Value = 38
$ws.Range("D64").Value = -1.046
$ws.Range("E64").Value = -2
$ws.Range("F64").Value = 38

# Row 68: U.S. Total
$ws.Range("B68").Value = -5933
$ws.Range("C68").Value = -4811
$ws.Range("D68").Value = 0.233
$ws.Range("E68").Value = -4972
$ws.Range("F68").Value = -3895
$ws.Range("G68").Value = -961
$ws.Range("H68").Value = -916
